$p = $ppt.ActivePresentation

$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(1).TextFrame.TextRange
$tr1.Text = "placeholder"
$tr1.Text = "Example numbering MWE"

$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(1).TextFrame.TextRange
$tr2.Text = "placeholder"
$tr2.Text = "A second slide"
